$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 67741.01150000001
$ws.Range("B2").Value = 55000
$ws.Range("C2").Value = 25000
$ws.Range("D2").Value = 20000
$ws.Range("E2").Value = 10000
$ws.Range("F2").Value = 10000
$ws.Range("G2").Value = 10000
$ws.Range("H2").Value = 7324.023000000008

$ws.Range("A3").Value = 64760.858
$ws.Range("H3").Value = 37050.716

$ws.Range("A4").Value = 62101.3
$ws.Range("H4").Value = 34469.60000000001

$ws.Range("A5").Value = 61724.802
$ws.Range("H5").Value = 34175.60400000001

$ws.Range("A6").Value = 63342.0075
$ws.Range("H6").Value = 35757.015

$ws.Range("A7").Value = 67452.89449999999
$ws.Range("H7").Value = 39572.789

$ws.Range("A8").Value = 66444.84849999999
$ws.Range("H8").Value = 37531.69699999999

$ws.Range("A9").Value = 77477.15949999999
$ws.Range("H9").Value = 47769.31899999999

$ws.Range("A10").Value = 93062.8
$ws.Range("B10").Value = 55000
$ws.Range("G10").Value = 10000
$ws.Range("H10").Value = 28573.60000000001

$ws.Range("A11").Value = 97938.8115
$ws.Range("H11").Value = 18620.62299999999

$ws.Range("A12").Value = 100816.808
$ws.Range("H12").Value = 20476.61600000001

$ws.Range("A13").Value = 99591.12850000001
$ws.Range("H13").Value = 19001.25700000001

$ws.Range("A14").Value = 103078.6315
$ws.Range("H14").Value = 22719.26300000001

$ws.Range("A15").Value = 103456.56
$ws.Range("H15").Value = 22945.12

$ws.Range("A16").Value = 106139.5625
$ws.Range("H16").Value = 25473.125

$ws.Range("A17").Value = 98840.583
$ws.Range("H17").Value = 17715.166

$ws.Range("A18").Value = 94650.12700000001
$ws.Range("H18").Value = 13243.25400000002

$ws.Range("A19").Value = 91921.8355
$ws.Range("H19").Value = 10758.671

$ws.Range("A20").Value = 91677.2825
$ws.Range("H20").Value = 11017.565

$ws.Range("A21").Value = 78811.4135
$ws.Range("H21").Value = 46.8269999999975

$ws.Range("A22").Value = 80527.792
$ws.Range("H22").Value = 4021.584000000003

$ws.Range("A23").Value = 59160.206
$ws.Range("H23").Value = -13942.588

$ws.Range("A24").Value = 65087.887
$ws.Range("H24").Value = -6220.225999999995

$ws.Range("H25").Value = 855.6514999999927
